$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => (new D value or $null if unchanged, new E value)
$updates = @(
    @{ Row = 2;  D = "27.156.23";      E = "  +0.52%  " },
    @{ Row = 3;  D = "1.829.87";       E = "  +0.30%  " },
    @{ Row = 4;  D = "1.010";          E = "  +0.45%  " },
    @{ Row = 5;  D = "313.02";         E = "  +0.52%  " },
    @{ Row = 6;  D = $null;            E = "  +0.38%  " },
    @{ Row = 7;  D = "0.4706";         E = "  +0.21%  " },
    @{ Row = 8;  D = "0.3670";         E = "  +0.07%  " },
    @{ Row = 9;  D = "0.07398";        E = "  +0.53%  " },
    @{ Row = 10; D = "0.8814";         E = "  +0.63%  " },
    @{ Row = 11; D = "20.31";          E = "  -0.01%  " },
    @{ Row = 12; D = "1.893.25";       E = "  +3.34%  " },
    @{ Row = 13; D = $null;            E = "  +4.71%  " },
    @{ Row = 14; D = "93.45";          E = "  +1.68%  " },
    @{ Row = 15; D = "5.381";          E = "  -0.98%  " },
    @{ Row = 16; D = "6.536";          E = "  +0.20%  " },
    @{ Row = 17; D = $null;            E = "  +0.04%  " },
    @{ Row = 18; D = "0.000008732";    E = "  -0.19%  " },
    @{ Row = 19; D = $null;            E = "  +0.37%  " },
    @{ Row = 20; D = "27.595.52";      E = "  +2.10%  " },
    @{ Row = 21; D = "14.64";          E = "  -0.43%  " },
    @{ Row = 22; D = $null;            E = "  -0.74%  " },
    @{ Row = 23; D = "10.63";          E = "  -0.05%  " },
    @{ Row = 24; D = "2.087.10";       E = "  +1.67%  " },
    @{ Row = 25; D = "1.881";          E = "  -0.65%  " },
    @{ Row = 26; D = "151.11";         E = "  -0.19%  " },
    @{ Row = 27; D = "18.52";          E = "  +0.41%  " },
    @{ Row = 28; D = "2.136";          E = "  -0.55%  " },
    @{ Row = 29; D = "5.181";          E = "  -1.17%  " },
    @{ Row = 30; D = "116.63";         E = "  -0.06%  " },
    @{ Row = 31; D = "0.08935";        E = "  +0.46%  " },
    @{ Row = 32; D = "0.7447";         E = "  -1.47%  " },
    @{ Row = 33; D = "1.165";          E = "  +0.26%  " },
    @{ Row = 34; D = "4.520";          E = "  +0.09%  " },
    @{ Row = 35; D = "2.941";          E = "  +0.30%  " },
    @{ Row = 36; D = "1.009";          E = "  +0.43%  " },
    @{ Row = 37; D = "2.563";          E = "  +7.54%  " },
    @{ Row = 38; D = $null;            E = "  -0.61%  " },
    @{ Row = 39; D = "0.05309";        E = "  -0.13%  " },
    @{ Row = 40; D = "0.01939";        E = "  -0.62%  " },
    @{ Row = 41; D = "7.340";          E = "  +1.57%  " },
    @{ Row = 42; D = $null;            E = "  -1.51%  " },
    @{ Row = 43; D = "0.5269";         E = "  -0.77%  " },
    @{ Row = 44; D = $null;            E = "  -0.72%  " },
    @{ Row = 45; D = "8.390";          E = "  -1.23%  " },
    @{ Row = 46; D = "0.4908";         E = "  +0.03%  " },
    @{ Row = 47; D = "10.42";          E = "  -0.14%  " },
    @{ Row = 48; D = "1.009";          E = "  +0.42%  " },
    @{ Row = 49; D = "104.44";         E = "  +1.23%  " },
    @{ Row = 50; D = "1.654";          E = "  -0.77%  " },
    @{ Row = 51; D = "0.06279";        E = "  -0.29%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        # Leading apostrophe forces text entry (matches the source data's
        # text-typed price strings like "27.156.23" / "1.010") without
        # Excel re-interpreting them as numbers. Resetting the style back
        # to Normal afterwards drops the quote-prefix formatting flag that
        # the apostrophe entry leaves behind, keeping cell styling as-is.
        $ws.Range("D$r").Value = "'" + $u.D
        $ws.Range("D$r").Style = "Normal"
    }
    $ws.Range("E$r").Value = $u.E
}
